$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 549 (shifts existing rows 549-561 down to 550-562)
$ws.Rows.Item(549).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A549").Value = 4
$ws.Range("B549").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C549").Value = "Los Lagos"
$ws.Range("D549").Value = 45239
$ws.Range("E549").Value = 10
$ws.Range("F549").Value = 100112045
$ws.Range("G549").Value = "Zapallo"
$ws.Range("H549").Value = "Paine"
$ws.Range("I549").Value = "1a (guarda)"
$ws.Range("J549").Value = 500
$ws.Range("K549").Value = 1200
$ws.Range("L549").Value = 1200
$ws.Range("M549").Value = 1200
$ws.Range("N549").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O549").Value = "Región de O'Higgins"
$ws.Range("P549").Value = 1200
$ws.Range("Q549").Value = 1
$ws.Range("R549").Value = "Hortaliza"
